$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns F, G, H, J with header/data pairs (dependent dropdown lists)
$ws.Range("G1").Value = "INC"
$ws.Range("G2").Value = "Hold"
$ws.Range("G3").Value = "Progress"

$ws.Range("F1").Value = "COM"
$ws.Range("F2").Value = "Completed"

$ws.Range("H1").Value = "INC"
$ws.Range("H2").Value = "COM"

$ws.Range("J1").Value = "ABC"
$ws.Range("J2").Value = "DEF"

# New sample data block further down the sheet
$ws.Range("A8").Value = 100
$ws.Range("A9").Value = 22
$ws.Range("B9").Value = "ABC"
$ws.Range("A10").Value = 33
$ws.Range("B10").Value = "DEF"
$ws.Range("A11").Value = 100
$ws.Range("A12").Value = 22
$ws.Range("A13").Value = 100

# Column widths for E and F (nearest settable values to the target 12.5703125 / 11.7109375)
$ws.Range("E1").ColumnWidth = 11.65
$ws.Range("F1").ColumnWidth = 10.8

# Update the selected cell shown in the saved view
$ws.Range("F15").Select()
